# Update the cached "datetimeFigureOut" date field text from 2/23/22 to
# 8/1/22 everywhere it appears: once on the Slide Master and once on each
# of the 11 slide layouts (each has a "Date Placeholder" shape holding the
# cached field text).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "2/23/22") {
                $shp.TextFrame.TextRange.Text = "8/1/22"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}
